{"js": "// Update the two-digit / one-digit division problems in the table.\n// Each \"NN\u00f7N=\" expression is replaced with its new value, matching the\n// commit's regenerated worksheet numbers.\nconst replacements = [\n  [\"83\u00f78=\", \"34\u00f73=\"],\n  [\"35\u00f75=\", \"86\u00f75=\"],\n  [\"28\u00f72=\", \"50\u00f72=\"],\n  [\"93\u00f73=\", \"10\u00f79=\"],\n  [\"41\u00f78=\", \"85\u00f75=\"],\n  [\"86\u00f73=\", \"45\u00f73=\"],\n  [\"28\u00f79=\", \"21\u00f74=\"],\n  [\"60\u00f74=\", \"27\u00f73=\"],\n  [\"22\u00f78=\", \"75\u00f73=\"],\n  [\"96\u00f79=\", \"53\u00f74=\"],\n  [\"17\u00f77=\", \"49\u00f75=\"],\n  [\"11\u00f74=\", \"39\u00f73=\"],\n  [\"58\u00f77=\", \"93\u00f77=\"],\n  [\"44\u00f75=\", \"40\u00f78=\"],\n  [\"76\u00f76=\", \"49\u00f77=\"],\n  [\"20\u00f78=\", \"89\u00f78=\"],\n  [\"98\u00f73=\", \"98\u00f74=\"],\n  [\"89\u00f77=\", \"85\u00f73=\"],\n  [\"68\u00f74=\", \"36\u00f74=\"],\n  [\"69\u00f74=\", \"61\u00f79=\"],\n  [\"67\u00f75=\", \"28\u00f78=\"],\n  [\"91\u00f76=\", \"54\u00f72=\"],\n  [\"34\u00f79=\", \"45\u00f79=\"],\n  [\"41\u00f76=\", \"56\u00f78=\"],\n  [\"80\u00f77=\", \"60\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit / one-digit division problems in the table.\n# Each \"NN\u00f7N=\" expression is replaced with its new value, matching the\n# commit's regenerated worksheet numbers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"83\u00f78=\"; New = \"34\u00f73=\" },\n    @{ Old = \"35\u00f75=\"; New = \"86\u00f75=\" },\n    @{ Old = \"28\u00f72=\"; New = \"50\u00f72=\" },\n    @{ Old = \"93\u00f73=\"; New = \"10\u00f79=\" },\n    @{ Old = \"41\u00f78=\"; New = \"85\u00f75=\" },\n    @{ Old = \"86\u00f73=\"; New = \"45\u00f73=\" },\n    @{ Old = \"28\u00f79=\"; New = \"21\u00f74=\" },\n    @{ Old = \"60\u00f74=\"; New = \"27\u00f73=\" },\n    @{ Old = \"22\u00f78=\"; New = \"75\u00f73=\" },\n    @{ Old = \"96\u00f79=\"; New = \"53\u00f74=\" },\n    @{ Old = \"17\u00f77=\"; New = \"49\u00f75=\" },\n    @{ Old = \"11\u00f74=\"; New = \"39\u00f73=\" },\n    @{ Old = \"58\u00f77=\"; New = \"93\u00f77=\" },\n    @{ Old = \"44\u00f75=\"; New = \"40\u00f78=\" },\n    @{ Old = \"76\u00f76=\"; New = \"49\u00f77=\" },\n    @{ Old = \"20\u00f78=\"; New = \"89\u00f78=\" },\n    @{ Old = \"98\u00f73=\"; New = \"98\u00f74=\" },\n    @{ Old = \"89\u00f77=\"; New = \"85\u00f73=\" },\n    @{ Old = \"68\u00f74=\"; New = \"36\u00f74=\" },\n    @{ Old = \"69\u00f74=\"; New = \"61\u00f79=\" },\n    @{ Old = \"67\u00f75=\"; New = \"28\u00f78=\" },\n    @{ Old = \"91\u00f76=\"; New = \"54\u00f72=\" },\n    @{ Old = \"34\u00f79=\"; New = \"45\u00f79=\" },\n    @{ Old = \"41\u00f76=\"; New = \"56\u00f78=\" },\n    @{ Old = \"80\u00f77=\"; New = \"60\u00f72=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
